$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).RowHeight = 20.25
$ws.Rows.Item(18).RowHeight = 19.5
$ws.Rows.Item(19).RowHeight = 19.5

$wb.Save()
